$wb = $excel.ActiveWorkbook

# Sheet 1: "Obras en general" - currently A1:X1, add Y1 (WhatsApp Profesional) and Z1 (WhatsApp Tramitador)
$ws1 = $wb.Worksheets.Item("Obras en general")
$ws1.Range("X1").Copy()
$ws1.Range("Y1:Z1").PasteSpecial(-4122)  # xlPasteFormats - copy the header style (bold) to the new cells
$excel.CutCopyMode = $false
$ws1.Range("Y1").Value = "WhatsApp Profesional"
$ws1.Range("Z1").Value = "WhatsApp Tramitador"

# Sheet 2: "Informes técnicos" - currently A1:O1, add P1 (WhatsApp Profesional) and Q1 (WhatsApp Tramitador)
$ws2 = $wb.Worksheets.Item("Informes técnicos")
$ws2.Range("O1").Copy()
$ws2.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats - copy the header style (bold) to the new cells
$excel.CutCopyMode = $false
$ws2.Range("P1").Value = "WhatsApp Profesional"
$ws2.Range("Q1").Value = "WhatsApp Tramitador"
